# Aula 63 - Nomeando os componentes do HTML
# Adds a new row (69) to the "anotacoes" sheet describing lesson 63, copying
# the formatting used by the previous entries, and moves the selection/
# viewport down to the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (fonts, fills, wrap text, etc.) from the previous
# entry's row (68) onto the new row (69) before writing the new content so
# that B69/C69 (style index 5) and D69/E69 (style index 1) match the rest
# of the table.
$ws.Range("B68:E68").Copy()
$ws.Range("B69:E69").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New entry: Aula 63 - Nomeando os componentes do HTML
$ws.Range("B69").Value = 63
$ws.Range("C69").Value = "11. Validação Back-End"
$ws.Range("D69").Value = "63. Nomeando os componentes do HTML"
$ws.Range("E69").Value = "foi abordado como nomear os componentes HTML através do arquivo messages.properties. Uma boa forma de centralizar os titulos de pagina HTML dos componentes, cabeçalhos de tabelas, etc"

# The row holds a single wrapped paragraph like row 68, so match its height.
$ws.Rows.Item(69).RowHeight = 30

# Move the view / selection down to the new last row, like the author did.
$ws.Range("D70").Select()
